$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a literal text value into a cell without altering its style.
# Plain Range.Value assignment lets the COM layer "smart type" numeric-looking
# strings (e.g. "45.60" -> 45.6, "6.08" -> 6.0800000000000001 binary noise, and
# a fresh cell style). Routing the literal through a temporary ="..." text
# formula and then Copy / PasteSpecial(values) bakes down to a plain shared-
# string cell with the original numFmt/style untouched - matching how the
# source data (t="inlineStr") is preserved as text end to end.
function Set-LiteralText($ws, $addr, $text) {
    $escaped = $text.Replace('"', '""')
    $cell = $ws.Range($addr)
    $cell.Formula = '="' + $escaped + '"'
    $cell.Copy()
    $cell.PasteSpecial(-4163)
}

Set-LiteralText $ws 'D2' '38.965.05'
Set-LiteralText $ws 'E2' '  -4.16%  '
Set-LiteralText $ws 'D3' '2.224.73'
Set-LiteralText $ws 'E3' '  -6.47%  '
Set-LiteralText $ws 'E4' '  +0.05%  '
Set-LiteralText $ws 'D5' '295.86'
Set-LiteralText $ws 'E5' '  -5.23%  '
Set-LiteralText $ws 'D6' '79.43'
Set-LiteralText $ws 'E6' '  -9.01%  '
Set-LiteralText $ws 'E7' '  -4.16%  '
Set-LiteralText $ws 'D9' '0.458'
Set-LiteralText $ws 'E9' '  -6.65%  '
Set-LiteralText $ws 'D10' '0.0768'
Set-LiteralText $ws 'E10' '  -6.47%  '
Set-LiteralText $ws 'D11' '27.65'
Set-LiteralText $ws 'E11' '  -10.63%  '
Set-LiteralText $ws 'D12' '45.60'
Set-LiteralText $ws 'E12' '  -13.96%  '
Set-LiteralText $ws 'E13' '  -1.22%  '
Set-LiteralText $ws 'D14' '2.573.02'
Set-LiteralText $ws 'E14' '  -6.22%  '
Set-LiteralText $ws 'D15' '6.08'
Set-LiteralText $ws 'E15' '  -7.85%  '
Set-LiteralText $ws 'D16' '14.00'
Set-LiteralText $ws 'E16' '  -6.73%  '
Set-LiteralText $ws 'D17' '2.241.58'
Set-LiteralText $ws 'E17' '  -4.50%  '
Set-LiteralText $ws 'D18' '0.712'
Set-LiteralText $ws 'E18' '  -5.95%  '
Set-LiteralText $ws 'D19' '38.905.77'
Set-LiteralText $ws 'E19' '  -4.06%  '
Set-LiteralText $ws 'D20' '0.0₃0855'
Set-LiteralText $ws 'E20' '  -6.02%  '
Set-LiteralText $ws 'D21' '5.71'
Set-LiteralText $ws 'E21' '  -7.34%  '
Set-LiteralText $ws 'D22' '64.83'
Set-LiteralText $ws 'E22' '  -5.98%  '
Set-LiteralText $ws 'E23' '  -9.38%  '
Set-LiteralText $ws 'D24' '224.00'
Set-LiteralText $ws 'E24' '  -4.66%  '
Set-LiteralText $ws 'E25' '  -0.13%  '
Set-LiteralText $ws 'D26' '2.37'
Set-LiteralText $ws 'E26' '  -10.08%  '
Set-LiteralText $ws 'D27' '1.72'
Set-LiteralText $ws 'E27' '  -5.48%  '
Set-LiteralText $ws 'D28' '22.13'
Set-LiteralText $ws 'E28' '  -6.05%  '
Set-LiteralText $ws 'D29' '2.14'
Set-LiteralText $ws 'E29' '  -2.93%  '
Set-LiteralText $ws 'D30' '8.87'
Set-LiteralText $ws 'E30' '  -5.34%  '
Set-LiteralText $ws 'D31' '149.34'
Set-LiteralText $ws 'E31' '  -3.93%  '
Set-LiteralText $ws 'D32' '30.93'
Set-LiteralText $ws 'E32' '  -8.28%  '
Set-LiteralText $ws 'E33' '  -0.08%  '
Set-LiteralText $ws 'E34' '  -8.82%  '
Set-LiteralText $ws 'E35' '  -3.12%  '
Set-LiteralText $ws 'E36' '  -6.34%  '
Set-LiteralText $ws 'E37' '  -4.23%  '
Set-LiteralText $ws 'E38' '  -5.64%  '
Set-LiteralText $ws 'E39' '  -3.77%  '
Set-LiteralText $ws 'B40' 'ARBITRUM'
Set-LiteralText $ws 'C40' 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-LiteralText $ws 'D40' '1.60'
Set-LiteralText $ws 'E40' '  -7.62%  '
Set-LiteralText $ws 'B41' 'Celestia'
Set-LiteralText $ws 'C41' 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
Set-LiteralText $ws 'D41' '14.43'
Set-LiteralText $ws 'E41' '  -9.63%  '
Set-LiteralText $ws 'E42' '  -4.72%  '
Set-LiteralText $ws 'D43' '1.901.23'
Set-LiteralText $ws 'E43' '  -3.11%  '
Set-LiteralText $ws 'D44' '2.07'
Set-LiteralText $ws 'E44' '  -8.88%  '
Set-LiteralText $ws 'E45' '  -6.14%  '
Set-LiteralText $ws 'D46' '16.21'
Set-LiteralText $ws 'E46' '  -8.10%  '
Set-LiteralText $ws 'D47' '8.96'
Set-LiteralText $ws 'E47' '  -3.59%  '
Set-LiteralText $ws 'D48' '2.50'
Set-LiteralText $ws 'E48' '  -10.63%  '
Set-LiteralText $ws 'D49' '2.442.69'
Set-LiteralText $ws 'E49' '  -6.21%  '
Set-LiteralText $ws 'D50' '87.24'
Set-LiteralText $ws 'E50' '  -6.65%  '
Set-LiteralText $ws 'D51' '67.40'
Set-LiteralText $ws 'E51' '  -7.27%  '

$excel.CutCopyMode = $false

